# 3.1b README.docx — correct a mistake:
#  - the paragraph "0x080003AC B08A      SUB      sp,sp,#0x28" and the
#    following paragraph "Затем адрес переменной в данной функции
#    переносится в другой регистр." (right before "Сейчас начинается
#    действие функции strToInt.") are emptied out.
#  - the stray "_GoBack" bookmark that sat between ". Пропускаем
#    стековый кадр и в" and "ыгружаем из стека " is removed from there
#    and re-created at the very start of the "Сейчас начинается
#    действие функции strToInt." paragraph instead.

$d = $word.ActiveDocument

# --- locate the two paragraphs to clear -----------------------------
# Find the unique spot where the asm-dump line is immediately followed
# (across the paragraph break) by "Затем адрес" — this disambiguates
# it from the earlier, identical-looking asm-dump line a few
# paragraphs above that must stay untouched.
$needle = "0x080003AC B08A      SUB      sp,sp,#0x28" + [char]13 + "Затем адрес"
$hit = $d.Content
$found = $hit.Find.Execute($needle, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "could not locate the '0x080003AC ... / Затем адрес ...' paragraphs"
}

$anchor = $d.Range($hit.Start, $hit.Start)
$paraAsm = $anchor.Paragraphs(1)        # "0x080003AC ... sp,#0x28"
$paraAddr = $paraAsm.Next()             # "Затем адрес ..."
$paraNext = $paraAddr.Next()            # "Сейчас начинается ..."

# Clear the later paragraph first so the earlier paragraph's
# position in the document is never disturbed by this edit.
$paraAddr.Range.InsertXML("<w:p><w:pPr></w:pPr></w:p>")
$paraAsm.Range.InsertXML("<w:p><w:pPr></w:pPr></w:p>")

# --- move the _GoBack bookmark --------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}
$bmPoint = $paraNext.Range.Start
$bmRange = $d.Range($bmPoint, $bmPoint)
$d.Bookmarks.Add("_GoBack", $bmRange)
